# Apply odds updates to Sheet1 for the 2025-01-27 FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 3.6
# Row 3
$ws.Range("AP3").Value = 1.89
$ws.Range("AQ3").Value = 2.01
# Row 4
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3.2
# Row 5
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
# Row 6
$ws.Range("AP6").Value = 1.88
$ws.Range("AQ6").Value = 1.98
# Row 7
$ws.Range("AR7").Value = 3.5
$ws.Range("AS7").Value = 1.29
# Row 8
$ws.Range("G8").Value = 1.39
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 1.88
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 7.8
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 6.7
$ws.Range("O8").Value = 1.35
$ws.Range("P8").Value = 2.95
$ws.Range("Q8").Value = 2.05
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 3.45
$ws.Range("T8").Value = 1.26
$ws.Range("X8").Value = 1.53
$ws.Range("AA8").Value = 8.75
$ws.Range("AB8").Value = 8.5
$ws.Range("AE8").Value = 6.7
$ws.Range("AF8").Value = 8.25
$ws.Range("AG8").Value = 25
$ws.Range("AJ8").Value = 18.5
$ws.Range("AL8").Value = 28
$ws.Range("AM8").Value = 300
$ws.Range("AO8").Value = 120
# Row 9
$ws.Range("G9").Value = 1.9
$ws.Range("H9").Value = 2.85
$ws.Range("J9").Value = 2.57
$ws.Range("K9").Value = 1.85
$ws.Range("L9").Value = 5.5
$ws.Range("N9").Value = 4.85
$ws.Range("P9").Value = 2.2
$ws.Range("Q9").Value = 2.75
$ws.Range("V9").Value = 2.18
$ws.Range("W9").Value = 2.32
$ws.Range("Y9").Value = 4.75
$ws.Range("Z9").Value = 7.3
$ws.Range("AB9").Value = 16
$ws.Range("AC9").Value = 20
$ws.Range("AE9").Value = 4.85
$ws.Range("AF9").Value = 6
$ws.Range("AJ9").Value = 9
$ws.Range("AL9").Value = 18
$ws.Range("AM9").Value = 110
$ws.Range("AN9").Value = 75
# Row 10
$ws.Range("G10").Value = 4.55
$ws.Range("I10").Value = 2.15
$ws.Range("J10").Value = 5.2
$ws.Range("L10").Value = 2.92
$ws.Range("X10").Value = 1.53
$ws.Range("AA10").Value = 16
$ws.Range("AC10").Value = 65
$ws.Range("AL10").Value = 9.75
$ws.Range("AM10").Value = 22
$ws.Range("AN10").Value = 25
$ws.Range("AO10").Value = 55
# Row 11
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 2.5
$ws.Range("Q11").Value = 2.63
$ws.Range("R11").Value = 1.5
$ws.Range("S11").Value = 5
$ws.Range("T11").Value = 1.17
$ws.Range("AP11").Value = 2
$ws.Range("AQ11").Value = 1.85
# Row 12
$ws.Range("G12").Value = 2.9
$ws.Range("I12").Value = 2.38
$ws.Range("J12").Value = 3.75
$ws.Range("L12").Value = 3.25
$ws.Range("O12").Value = 1.4
$ws.Range("P12").Value = 2.75
$ws.Range("Q12").Value = 2.3
$ws.Range("R12").Value = 1.6
$ws.Range("S12").Value = 4.33
$ws.Range("T12").Value = 1.2
$ws.Range("W12").Value = 1.91
$ws.Range("X12").Value = 1.8
$ws.Range("Y12").Value = 8.5
$ws.Range("Z12").Value = 15
$ws.Range("AA12").Value = 12
$ws.Range("AB12").Value = 34
$ws.Range("AC12").Value = 26
$ws.Range("AD12").Value = 41
$ws.Range("AE12").Value = 7.5
$ws.Range("AF12").Value = 6
$ws.Range("AI12").Value = 351
$ws.Range("AJ12").Value = 7.5
$ws.Range("AK12").Value = 11
$ws.Range("AL12").Value = 10
$ws.Range("AM12").Value = 23
$ws.Range("AN12").Value = 21
# Row 13
$ws.Range("G13").Value = 1.55
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 5.25
$ws.Range("J13").Value = 2.1
$ws.Range("L13").Value = 5.5
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 4.33
$ws.Range("Q13").Value = 1.67
$ws.Range("R13").Value = 2.15
$ws.Range("W13").Value = 1.8
$ws.Range("X13").Value = 1.91
$ws.Range("Y13").Value = 7.5
$ws.Range("Z13").Value = 8
$ws.Range("AC13").Value = 12
$ws.Range("AF13").Value = 8
$ws.Range("AG13").Value = 17
$ws.Range("AH13").Value = 51
# Row 15
$ws.Range("G15").Value = 2.4
$ws.Range("I15").Value = 2.5
$ws.Range("L15").Value = 3.1
$ws.Range("Y15").Value = 12
$ws.Range("Z15").Value = 15
$ws.Range("AA15").Value = 10
$ws.Range("AB15").Value = 26
$ws.Range("AC15").Value = 19
$ws.Range("AD15").Value = 23
$ws.Range("AM15").Value = 26
# Row 16
$ws.Range("G16").Value = 1.57
$ws.Range("H16").Value = 4.2
$ws.Range("L16").Value = 5
$ws.Range("W16").Value = 1.67
$ws.Range("X16").Value = 2.1
$ws.Range("AF16").Value = 8
$ws.Range("AR16").Value = 2.05
$ws.Range("AS16").Value = 1.8
# Row 17
$ws.Range("G17").Value = 2.25
$ws.Range("I17").Value = 2.75
$ws.Range("J17").Value = 2.75
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 1.73
$ws.Range("AC17").Value = 15
$ws.Range("AK17").Value = 19
$ws.Range("AL17").Value = 11
$ws.Range("AN17").Value = 19
# Row 18
$ws.Range("I18").Value = 2.4
$ws.Range("Q18").Value = 2.35
$ws.Range("R18").Value = 1.57
$ws.Range("U18").Value = 1.53
$ws.Range("V18").Value = 2.38
$ws.Range("Y18").Value = 8
$ws.Range("AG18").Value = 17
$ws.Range("AI18").Value = 401
# Row 20
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 4.5
$ws.Range("I20").Value = 1.55
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 2.38
$ws.Range("L20").Value = 2.05
$ws.Range("M20").Value = 1.03
$ws.Range("N20").Value = 10.5
$ws.Range("O20").Value = 1.17
$ws.Range("P20").Value = 4.5
$ws.Range("Q20").Value = 1.57
$ws.Range("R20").Value = 2.35
$ws.Range("S20").Value = 2.38
$ws.Range("T20").Value = 1.53
$ws.Range("U20").Value = 1.3
$ws.Range("V20").Value = 3.4
$ws.Range("W20").Value = 1.73
$ws.Range("X20").Value = 2
$ws.Range("Y20").Value = 17
$ws.Range("AD20").Value = 41
$ws.Range("AE20").Value = 15
$ws.Range("AF20").Value = 9
$ws.Range("AH20").Value = 51
$ws.Range("AJ20").Value = 8.5
$ws.Range("AL20").Value = 8.5
$ws.Range("AN20").Value = 12
$ws.Range("AO20").Value = 23
$ws.Range("AR20").Value = 1.95
$ws.Range("AS20").Value = 1.85
# Row 21
$ws.Range("G21").Value = 1.91
$ws.Range("H21").Value = 3.9
$ws.Range("I21").Value = 3.5
$ws.Range("K21").Value = 2.2
$ws.Range("L21").Value = 4
$ws.Range("M21").Value = 1.04
$ws.Range("N21").Value = 9
$ws.Range("O21").Value = 1.25
$ws.Range("P21").Value = 3.75
$ws.Range("Q21").Value = 1.83
$ws.Range("R21").Value = 1.98
$ws.Range("S21").Value = 3
$ws.Range("T21").Value = 1.36
$ws.Range("W21").Value = 1.8
$ws.Range("X21").Value = 1.91
$ws.Range("Y21").Value = 7.5
$ws.Range("Z21").Value = 9
$ws.Range("AC21").Value = 15
$ws.Range("AD21").Value = 26
$ws.Range("AE21").Value = 11
$ws.Range("AI21").Value = 600
$ws.Range("AJ21").Value = 11
$ws.Range("AO21").Value = 34
# Row 22
$ws.Range("U22").Value = 1.29
$ws.Range("V22").Value = 3.5
$ws.Range("W22").Value = 1.53
$ws.Range("X22").Value = 2.38
$ws.Range("Y22").Value = 15
$ws.Range("AK22").Value = 11
# Row 24
$ws.Range("G24").Value = 1.67
$ws.Range("H24").Value = 4.2
$ws.Range("I24").Value = 4.5
$ws.Range("O24").Value = 1.18
$ws.Range("P24").Value = 4.5
$ws.Range("Q24").Value = 1.62
$ws.Range("R24").Value = 2.25
$ws.Range("S24").Value = 2.5
$ws.Range("T24").Value = 1.5
$ws.Range("U24").Value = 1.3
$ws.Range("V24").Value = 3.4
$ws.Range("Z24").Value = 9
$ws.Range("AC24").Value = 12
